# Generate Report for Handoff
# Update Priority from "low" to "ht" and refresh the Latest Handoff Datetime
# for the still-pending rows (4-7) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-01 08:39:12"
}

# de-de sheet: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> new timestamp
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-01 08:39:18"
}

# Overview sheet: rows 4-7, Latest HO Xliff Generate Date (G) shares the same
# underlying text as de-de's Latest Handoff Datetime, so it advances too.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-09-01 08:39:18"
}
